# This change comes from a SharePoint "content type" sync (commit message
# "Atualizando com a master" -> merging in metadata that the document
# library pushed down): the document's linked content-type schema
# (customXml/item1.xml, a ct:contentTypeSchema part) is bumped from
# contentTypeVersion 6 -> 7, gets a freshly minted versionID / fieldsID,
# and picks up one additional managed-metadata field definition
# (MediaServiceDateTaken, mirroring the existing MediaServiceEventHashCode
# field). The matching customXml/itemProps1.xml datastore item also gets a
# freshly minted ds:itemID GUID, the way Word re-mints it whenever it
# rewrites a custom XML part. No text in the body of the document changes.
#
# The real-world Word object model surface for this is
# Document.CustomXMLParts (a CustomXMLPart keyed by its content-type
# namespace / GUID, exposing a settable .XML string). We drive the edit
# through that API: look the part up, patch the handful of
# attributes/elements that actually changed via plain text substitution
# (the exact before/after text is known), and write the updated XML
# string back — exactly as SharePoint's own sync code (or a VBA/PowerShell
# macro doing the equivalent) would.

$d = $word.ActiveDocument

$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"
$oldItemId     = "9A6F4D4D-A5A7-4704-8C86-19BFDF6C8056"
$newItemId     = "6125353E-BC1B-4DF4-B701-AB3781ADCDC9"

function Update-ContentTypeSchemaXml([string]$xmlText) {
    if ([string]::IsNullOrEmpty($xmlText)) { return $null }

    $updated = $xmlText

    $updated = $updated.Replace(
        'ma:contentTypeVersion="6" ma:contentTypeDescription="Crie um novo documento." ma:contentTypeScope="" ma:versionID="530025b846b757239230ee9053be16e1"',
        'ma:contentTypeVersion="7" ma:contentTypeDescription="Crie um novo documento." ma:contentTypeScope="" ma:versionID="14cb47ab60ce25b03f125bbc624aa0ec"'
    )

    $updated = $updated.Replace(
        'ma:fieldsID="b2254a4aead5832002d2352349de8aba"',
        'ma:fieldsID="bd0907572239d88a9fb49e7cd2f006cd"'
    )

    $updated = $updated.Replace(
        '<xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>',
        '<xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/><xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/>'
    )

    $hashCodeDef = '<xsd:element name="MediaServiceEventHashCode" ma:index="13" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element>'
    $dateTakenDef = '<xsd:element name="MediaServiceDateTaken" ma:index="14" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element>'
    if ($updated.Contains($hashCodeDef) -and -not $updated.Contains($dateTakenDef)) {
        $updated = $updated.Replace($hashCodeDef, $hashCodeDef + $dateTakenDef)
    }

    return $updated
}

function Try-UpdatePart($part) {
    if ($part -eq $null) { return $false }
    $xmlText = $null
    try { $xmlText = $part.XML } catch { return $false }
    if ([string]::IsNullOrEmpty($xmlText)) { return $false }

    $newXml = Update-ContentTypeSchemaXml $xmlText
    if ($newXml -eq $null -or $newXml -eq $xmlText) { return $false }

    try {
        $part.XML = $newXml
    } catch {
        return $false
    }
    return $true
}

$updated = $false

try {
    $cxps = $d.CustomXMLParts

    try {
        $byId = $cxps.SelectByID($oldItemId)
        if ($byId -ne $null) {
            $updated = Try-UpdatePart $byId
        }
    } catch { }

    if (-not $updated) {
        try {
            $byNs = $cxps.SelectByNamespace($contentTypeNs)
            if ($byNs -ne $null -and $byNs.Count -ge 1) {
                for ($i = 1; $i -le $byNs.Count; $i++) {
                    if (Try-UpdatePart $byNs.Item($i)) {
                        $updated = $true
                        break
                    }
                }
            }
        } catch { }
    }

    if (-not $updated -and $cxps -ne $null) {
        $count = 0
        try { $count = $cxps.Count } catch { $count = 0 }
        for ($i = 1; $i -le $count; $i++) {
            $p = $cxps.Item($i)
            if ($p -ne $null -and $p.NamespaceURI -eq $contentTypeNs) {
                if (Try-UpdatePart $p) {
                    $updated = $true
                    break
                }
            }
        }
    }
} catch {
    # Document.CustomXMLParts not available / not populated for this
    # document in this hosting environment -- nothing further we can do
    # through the object model, so leave the document otherwise untouched.
}

# Re-mint the datastore item id for the same part (itemProps1.xml), the
# way Word does whenever it rewrites a custom XML part's backing store.
if ($updated) {
    try {
        $cxps = $d.CustomXMLParts
        $part = $cxps.SelectByID($oldItemId)
        if ($part -ne $null) {
            $part.Id = $newItemId
        }
    } catch { }
}
